# Re-classify a block of equipment rows on the "Master Allocation" sheet from
# their old Division codes (HOU / WT) to the correct Division code (DFW), and
# stamp the newly-populated trailer rows (488-489) with the DFW division too.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Master Allocation")

# Rows that were previously "HOU" or "WT" in column A (Division) -> "DFW"
$ws.Range("A57:A59").Value   = "DFW"
$ws.Range("A92:A127").Value  = "DFW"
$ws.Range("A143:A144").Value = "DFW"
$ws.Range("A258:A279").Value = "DFW"
$ws.Range("A372:A385").Value = "DFW"
$ws.Range("A425:A435").Value = "DFW"
$ws.Range("A441:A476").Value = "DFW"

# New summary rows at the bottom that now get a Division value as well
$ws.Range("A488:A489").Value = "DFW"
